$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 110.9189436435699

$ws.Range("A4").Value = 18816.2325
$ws.Range("B4").Value = 18832
$ws.Range("C4").Value = 18622
$ws.Range("D4").Value = 18830
$ws.Range("E4").Value = 18589
$ws.Range("F4").Value = 6412.6795
$ws.Range("G4").Value = 6491
$ws.Range("H4").Value = 6334
$ws.Range("I4").Value = 6464
$ws.Range("J4").Value = 6301
